$d = $word.ActiveDocument

# Update the date header
$d.Paragraphs.Item(1).Range.Find.Execute("2023-03-26 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-27 Monday", 2) | Out-Null

# Update each multiplication problem cell by exact (row, col) address
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "72×89="
$t.Cell(1, 2).Range.Text = "86×96="
$t.Cell(1, 3).Range.Text = "39×14="
$t.Cell(1, 4).Range.Text = "17×12="
$t.Cell(1, 5).Range.Text = "79×20="
$t.Cell(2, 1).Range.Text = "98×91="
$t.Cell(2, 2).Range.Text = "55×56="
$t.Cell(2, 3).Range.Text = "20×10="
$t.Cell(2, 4).Range.Text = "99×33="
$t.Cell(2, 5).Range.Text = "53×25="
$t.Cell(3, 1).Range.Text = "95×58="
$t.Cell(3, 2).Range.Text = "39×42="
$t.Cell(3, 3).Range.Text = "57×74="
$t.Cell(3, 4).Range.Text = "91×82="
$t.Cell(3, 5).Range.Text = "70×68="
$t.Cell(4, 1).Range.Text = "79×51="
$t.Cell(4, 2).Range.Text = "14×36="
$t.Cell(4, 3).Range.Text = "34×95="
$t.Cell(4, 4).Range.Text = "86×26="
$t.Cell(4, 5).Range.Text = "92×27="
$t.Cell(5, 1).Range.Text = "56×63="
$t.Cell(5, 2).Range.Text = "72×76="
$t.Cell(5, 3).Range.Text = "52×50="
$t.Cell(5, 4).Range.Text = "85×15="
$t.Cell(5, 5).Range.Text = "99×78="
$t.Cell(6, 1).Range.Text = "64×31="
$t.Cell(6, 2).Range.Text = "33×43="
$t.Cell(6, 3).Range.Text = "27×33="
$t.Cell(6, 4).Range.Text = "81×35="
$t.Cell(6, 5).Range.Text = "71×87="
$t.Cell(7, 1).Range.Text = "69×23="
$t.Cell(7, 2).Range.Text = "38×14="
$t.Cell(7, 3).Range.Text = "72×81="
$t.Cell(7, 4).Range.Text = "13×71="
$t.Cell(7, 5).Range.Text = "49×48="
$t.Cell(8, 1).Range.Text = "23×65="
$t.Cell(8, 2).Range.Text = "61×83="
$t.Cell(8, 3).Range.Text = "93×11="
$t.Cell(8, 4).Range.Text = "55×91="
$t.Cell(8, 5).Range.Text = "67×93="
$t.Cell(9, 1).Range.Text = "27×60="
$t.Cell(9, 2).Range.Text = "90×58="
$t.Cell(9, 3).Range.Text = "55×80="
$t.Cell(9, 4).Range.Text = "72×70="
$t.Cell(9, 5).Range.Text = "96×34="
$t.Cell(10, 1).Range.Text = "52×10="
$t.Cell(10, 2).Range.Text = "94×49="
$t.Cell(10, 3).Range.Text = "52×52="
$t.Cell(10, 4).Range.Text = "12×36="
$t.Cell(10, 5).Range.Text = "97×27="
$t.Cell(11, 1).Range.Text = "66×68="
$t.Cell(11, 2).Range.Text = "46×76="
$t.Cell(11, 3).Range.Text = "17×83="
$t.Cell(11, 4).Range.Text = "50×25="
$t.Cell(11, 5).Range.Text = "45×87="
$t.Cell(12, 1).Range.Text = "76×82="
$t.Cell(12, 2).Range.Text = "79×70="
$t.Cell(12, 3).Range.Text = "88×45="
$t.Cell(12, 4).Range.Text = "36×22="
$t.Cell(12, 5).Range.Text = "31×55="
$t.Cell(13, 1).Range.Text = "56×93="
$t.Cell(13, 2).Range.Text = "62×97="
$t.Cell(13, 3).Range.Text = "22×95="
$t.Cell(13, 4).Range.Text = "14×55="
$t.Cell(13, 5).Range.Text = "14×23="
$t.Cell(14, 1).Range.Text = "97×37="
$t.Cell(14, 2).Range.Text = "43×91="
$t.Cell(14, 3).Range.Text = "97×23="
$t.Cell(14, 4).Range.Text = "87×79="
$t.Cell(14, 5).Range.Text = "21×100="
$t.Cell(15, 1).Range.Text = "99×38="
$t.Cell(15, 2).Range.Text = "34×41="
$t.Cell(15, 3).Range.Text = "55×56="
$t.Cell(15, 4).Range.Text = "20×56="
$t.Cell(15, 5).Range.Text = "98×37="
$t.Cell(16, 1).Range.Text = "67×91="
$t.Cell(16, 2).Range.Text = "25×38="
$t.Cell(16, 3).Range.Text = "47×16="
$t.Cell(16, 4).Range.Text = "52×56="
$t.Cell(16, 5).Range.Text = "91×45="
$t.Cell(17, 1).Range.Text = "77×15="
$t.Cell(17, 2).Range.Text = "74×93="
$t.Cell(17, 3).Range.Text = "20×53="
$t.Cell(17, 4).Range.Text = "55×47="
$t.Cell(17, 5).Range.Text = "50×17="
$t.Cell(18, 1).Range.Text = "51×80="
$t.Cell(18, 2).Range.Text = "34×99="
$t.Cell(18, 3).Range.Text = "50×16="
$t.Cell(18, 4).Range.Text = "18×82="
$t.Cell(18, 5).Range.Text = "43×68="
$t.Cell(19, 1).Range.Text = "20×72="
$t.Cell(19, 2).Range.Text = "76×98="
$t.Cell(19, 3).Range.Text = "24×61="
$t.Cell(19, 4).Range.Text = "96×20="
$t.Cell(19, 5).Range.Text = "48×96="
$t.Cell(20, 1).Range.Text = "19×66="
$t.Cell(20, 2).Range.Text = "19×84="
$t.Cell(20, 3).Range.Text = "65×45="
$t.Cell(20, 4).Range.Text = "71×85="
$t.Cell(20, 5).Range.Text = "15×100="
